# Auto update ESO load data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row formatting: bold, thin box border, centered / top aligned ---
$header = $ws.Range("A1:AA1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108   # xlCenter
$header.VerticalAlignment = -4160     # xlTop
$header.Borders.LineStyle = 1         # xlContinuous
$header.Borders.Weight = 2            # xlThin

# --- Updated ESO load values (columns D:AA = hours 1-24, rows 2-7) ---
$newData = @{
    2 = @(5295,5047,4972,4941,4959,5122,5522,6000,6391,6548,6437,6224,6018,5848,5752,5815,5837,6120,6519,6509,6302,6021,5760,5419)
    3 = @(5028,4782,4708,4676,4694,4857,5255,5731,6120,6277,6176,5984,5798,5644,5558,5609,5627,5855,6176,6167,5963,5686,5429,5093)
    4 = @(4783,4498,4380,4394,4389,4412,4549,4926,5342,5650,5730,5616,5533,5332,5164,5152,5277,5529,5742,5625,5377,5112,4982,4749)
    5 = @(4713,4500,4402,4344,4318,4365,4473,4746,5112,5481,5675,5707,5683,5361,5199,5221,5382,5693,5946,5852,5637,5347,5156,4828)
    6 = @(4983,4728,4617,4628,4663,4797,5189,5667,6134,6381,6405,6283,6112,6061,5891,5837,5700,5926,6385,6361,6154,5865,5595,5226)
    7 = @(5248,5007,4934,4903,4922,5096,5523,6033,6450,6618,6511,6307,6109,5945,5853,5914,5935,6207,6592,6581,6365,6072,5801,5445)
}

foreach ($row in $newData.Keys) {
    $values = $newData[$row]
    $col = 4   # column D
    foreach ($v in $values) {
        $ws.Cells.Item($row, $col).Value = $v
        $col = $col + 1
    }
}
